# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# Summary of change (per the OOXML diff):
#   - "Body" sheet: the Request Body's "searchCriteria" object (rows 3-6) is
#     collapsed into a single schema-reference row pointing at
#     "interestDailyReport.211207Request".
#   - "200" sheet: the Response body detail rows (rows 3-11) are collapsed
#     into a single schema-reference row pointing at
#     "interestDailyReport.211207Response".
#   - "204" sheet: gains a new row 3 that references the
#     "interestDailyReport.211207Response" schema.
#   - "400" sheet: the detailed error rows (rows 3-6) are collapsed into a
#     single schema-reference row pointing at "errorResponse".
#   - "401", "403", "404", "429", "500" sheets: each gains a new row 3 that
#     references the "errorResponse1" schema.

$wb = $excel.ActiveWorkbook

function Set-SchemaRow($ws, $Section, $Name) {
    $ws.Range("A3").Value = $Section
    $ws.Range("B3").Value = $Name
    $ws.Range("C3").Value = ""
    $ws.Range("D3").Value = ""
    $ws.Range("E3").Value = "schema"
    $ws.Range("F3").Value = ""
    $ws.Range("G3").Value = $Name
    $ws.Range("H3").Value = ""
    $ws.Range("I3").Value = "Yes"
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = ""
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = ""
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = ""
}

# ---- "Body" sheet: collapse rows 3-6 (searchCriteria/settlementBIC/month) ----
# into a single row 3 referencing the interestDailyReport.211207Request schema.
$wsBody = $wb.Worksheets.Item("Body")
$wsBody.Rows.Item(4).Delete()
$wsBody.Rows.Item(4).Delete()
$wsBody.Rows.Item(4).Delete()
Set-SchemaRow $wsBody "body" "interestDailyReport.211207Request"

# ---- "200" sheet: collapse rows 3-11 (dateTime/settlementBIC/interestReports/
# dailyReport/date/interestRate/snapshotTime/eodPosition/dailyInterestAmount)
# into a single row 3 referencing the interestDailyReport.211207Response schema.
$ws200 = $wb.Worksheets.Item("200")
$ws200.Rows.Item(4).Delete()
$ws200.Rows.Item(4).Delete()
$ws200.Rows.Item(4).Delete()
$ws200.Rows.Item(4).Delete()
$ws200.Rows.Item(4).Delete()
$ws200.Rows.Item(4).Delete()
$ws200.Rows.Item(4).Delete()
$ws200.Rows.Item(4).Delete()
Set-SchemaRow $ws200 "content" "interestDailyReport.211207Response"

# ---- "204" sheet: add a new row 3 referencing interestDailyReport.211207Response.
$ws204 = $wb.Worksheets.Item("204")
Set-SchemaRow $ws204 "content" "interestDailyReport.211207Response"

# ---- "400" sheet: collapse rows 3-6 (dateTime/errorCode/errorCodeDescription/
# requestId) into a single row 3 referencing the errorResponse schema.
$ws400 = $wb.Worksheets.Item("400")
$ws400.Rows.Item(4).Delete()
$ws400.Rows.Item(4).Delete()
$ws400.Rows.Item(4).Delete()
Set-SchemaRow $ws400 "content" "errorResponse"

# ---- "401", "403", "404", "429", "500" sheets: each gains a new row 3
# referencing the errorResponse1 schema.
$ws401 = $wb.Worksheets.Item("401")
Set-SchemaRow $ws401 "content" "errorResponse1"

$ws403 = $wb.Worksheets.Item("403")
Set-SchemaRow $ws403 "content" "errorResponse1"

$ws404 = $wb.Worksheets.Item("404")
Set-SchemaRow $ws404 "content" "errorResponse1"

$ws429 = $wb.Worksheets.Item("429")
Set-SchemaRow $ws429 "content" "errorResponse1"

$ws500 = $wb.Worksheets.Item("500")
Set-SchemaRow $ws500 "content" "errorResponse1"
